$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 51422.348
$ws.Range("J17").Value = 48120.62
$ws.Range("L17").Value = 144361.86
$ws.Range("N17").Value = -144697.86
$ws.Range("H40").Value = 1771.5714
$ws.Range("I40").Value = 1801
$ws.Range("J40").Value = 1766.6666
$ws.Range("K40").Value = 1801
$ws.Range("L40").Value = 1766.6666
$ws.Range("M40").Value = -1626
$ws.Range("N40").Value = -2116.6666
$ws.Range("H58").Value = 2127.8125
$ws.Range("I58").Value = 361.875
$ws.Range("J58").Value = 2481
$ws.Range("K58").Value = 1085.625
$ws.Range("L58").Value = 7443
$ws.Range("M58").Value = -935.625
$ws.Range("N58").Value = -7743
$ws.Range("H62").Value = 2282.9524
$ws.Range("I62").Value = 1770.6154
$ws.Range("J62").Value = 3115.5
$ws.Range("K62").Value = 1770.6154
$ws.Range("L62").Value = 3115.5
$ws.Range("M62").Value = -1146.6154
$ws.Range("N62").Value = -4363.5
$ws.Range("H65").Value = 2282.9524
$ws.Range("I65").Value = 1770.6154
$ws.Range("J65").Value = 3115.5
$ws.Range("K65").Value = 8853.076999999999
$ws.Range("L65").Value = 15577.5
$ws.Range("M65").Value = -5733.076999999999
$ws.Range("N65").Value = -21817.5
$ws.Range("H82").Value = 7280
$ws.Range("I82").Value = 6800
$ws.Range("K82").Value = 20400
$ws.Range("M82").Value = -19994
$ws.Range("H85").Value = 7280
$ws.Range("I85").Value = 6800
$ws.Range("K85").Value = 20400
$ws.Range("M85").Value = -18996
$ws.Range("H115").Value = 1431
$ws.Range("I115").Value = 360.83334
$ws.Range("J115").Value = 2144.4443
$ws.Range("K115").Value = 1082.50002
$ws.Range("L115").Value = 6433.3329
$ws.Range("M115").Value = 484.4999800000001
$ws.Range("N115").Value = -9567.332900000001
$ws.Range("H116").Value = 3388.8125
$ws.Range("I116").Value = 2595.8
$ws.Range("J116").Value = 4710.5
$ws.Range("K116").Value = 2595.8
$ws.Range("L116").Value = 4710.5
$ws.Range("M116").Value = 846.1999999999998
$ws.Range("N116").Value = -11594.5
$ws.Range("H135").Value = 1262.0416
$ws.Range("I135").Value = 1059.6666
$ws.Range("J135").Value = 2678.6667
$ws.Range("K135").Value = 9536.999400000001
$ws.Range("L135").Value = 24108.0003
$ws.Range("M135").Value = -7001.999400000001
$ws.Range("N135").Value = -29178.0003
$ws.Range("H138").Value = 3042.3447
$ws.Range("I138").Value = 1710.8718
$ws.Range("J138").Value = 4124.1665
$ws.Range("K138").Value = 5132.6154
$ws.Range("L138").Value = 12372.4995
$ws.Range("M138").Value = 7.384600000000319
$ws.Range("N138").Value = -22652.4995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2207.5
$ws.Range("I61").Value = 1118.5769
$ws.Range("J61").Value = 3977
$ws.Range("K61").Value = 1118.5769
$ws.Range("L61").Value = 3977
$ws.Range("M61").Value = -906.5769
$ws.Range("N61").Value = -4401
$ws.Range("H74").Value = 874.3158
$ws.Range("I74").Value = 828.8
$ws.Range("J74").Value = 1045
$ws.Range("K74").Value = 828.8
$ws.Range("L74").Value = 1045
$ws.Range("M74").Value = 45.20000000000005
$ws.Range("N74").Value = -2793
$ws.Range("H77").Value = 874.3158
$ws.Range("I77").Value = 828.8
$ws.Range("J77").Value = 1045
$ws.Range("K77").Value = 4144
$ws.Range("L77").Value = 5225
$ws.Range("M77").Value = 224
$ws.Range("N77").Value = -13961
$ws.Range("H133").Value = 29792
$ws.Range("J133").Value = 29792
$ws.Range("L133").Value = 29792
$ws.Range("N133").Value = -34852
$ws.Range("H136").Value = 2207.5
$ws.Range("I136").Value = 1118.5769
$ws.Range("J136").Value = 3977
$ws.Range("K136").Value = 3355.7307
$ws.Range("L136").Value = 11931
$ws.Range("M136").Value = -805.7307000000001
$ws.Range("N136").Value = -17031

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 29842.857
$ws.Range("J115").Value = 29842.857
$ws.Range("L115").Value = 29842.857
$ws.Range("N115").Value = -32976.857

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 16596.416
$ws.Range("J74").Value = 16596.416
$ws.Range("L74").Value = 16596.416
$ws.Range("N74").Value = -18344.416
$ws.Range("H77").Value = 16596.416
$ws.Range("J77").Value = 16596.416
$ws.Range("L77").Value = 49789.24800000001
$ws.Range("N77").Value = -58525.24800000001
$ws.Range("H108").Value = 27800
$ws.Range("I108").Value = 15000
$ws.Range("J108").Value = 31000
$ws.Range("K108").Value = 15000
$ws.Range("L108").Value = 31000
$ws.Range("M108").Value = -11160
$ws.Range("N108").Value = -38680
$ws.Range("H122").Value = 2726.5
$ws.Range("I122").Value = 2288.7368
$ws.Range("J122").Value = 3914.7144
$ws.Range("K122").Value = 6866.2104
$ws.Range("L122").Value = 11744.1432
$ws.Range("M122").Value = -4416.2104
$ws.Range("N122").Value = -16644.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4923.273
$ws.Range("I3").Value = 5143.3335
$ws.Range("J3").Value = 4840.75
$ws.Range("K3").Value = 15430.0005
$ws.Range("L3").Value = 14522.25
$ws.Range("M3").Value = -15318.0005
$ws.Range("N3").Value = -14746.25
$ws.Range("H87").Value = 9053.467000000001
$ws.Range("I87").Value = 4544.6665
$ws.Range("K87").Value = 13633.9995
$ws.Range("M87").Value = -12385.9995
$ws.Range("H90").Value = 9053.467000000001
$ws.Range("I90").Value = 4544.6665
$ws.Range("K90").Value = 40901.9985
$ws.Range("M90").Value = -34661.9985
$ws.Range("H93").Value = 2349.4546
$ws.Range("J93").Value = 2327.2222
$ws.Range("L93").Value = 6981.6666
$ws.Range("N93").Value = -10725.6666
$ws.Range("H101").Value = 2918.182
$ws.Range("J101").Value = 2918.182
$ws.Range("L101").Value = 8754.545999999998
$ws.Range("N101").Value = -13622.546
$ws.Range("H120").Value = 18507.428
$ws.Range("H131").Value = 1885.4706
$ws.Range("I131").Value = 12000
$ws.Range("J131").Value = 1253.3125
$ws.Range("K131").Value = 36000
$ws.Range("L131").Value = 3759.9375
$ws.Range("M131").Value = -30960
$ws.Range("N131").Value = -13839.9375
$ws.Range("H138").Value = 1434.75
$ws.Range("I138").Value = 783.2222
$ws.Range("J138").Value = 3389.3333
$ws.Range("K138").Value = 2349.6666
$ws.Range("L138").Value = 10167.9999
$ws.Range("M138").Value = 2790.3334
$ws.Range("N138").Value = -20447.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 30014
$ws.Range("J26").Value = 30014
$ws.Range("L26").Value = 30014
$ws.Range("N26").Value = -30574
$ws.Range("H50").Value = 30014
$ws.Range("J50").Value = 30014
$ws.Range("L50").Value = 30014
$ws.Range("N50").Value = -31010

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 76924800
$ws.Range("I22").Value = 200000670
$ws.Range("J22").Value = 2372.75
$ws.Range("K22").Value = 200000670
$ws.Range("L22").Value = 2372.75
$ws.Range("M22").Value = -200000375
$ws.Range("N22").Value = -2962.75
$ws.Range("H27").Value = 76924800
$ws.Range("I27").Value = 200000670
$ws.Range("J27").Value = 2372.75
$ws.Range("K27").Value = 200000670
$ws.Range("L27").Value = 2372.75
$ws.Range("M27").Value = -200000563
$ws.Range("N27").Value = -2586.75
$ws.Range("H40").Value = 2650.4
$ws.Range("I40").Value = 1602
$ws.Range("J40").Value = 2912.5
$ws.Range("K40").Value = 1602
$ws.Range("L40").Value = 2912.5
$ws.Range("M40").Value = -1466
$ws.Range("N40").Value = -3184.5
$ws.Range("H93").Value = 5825
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 7266.6665
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 7266.6665
$ws.Range("M93").Value = -252
$ws.Range("N93").Value = -9762.666499999999
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
$ws.Range("H109").Value = 24320
$ws.Range("J109").Value = 24320
$ws.Range("L109").Value = 24320
$ws.Range("N109").Value = -27094
$ws.Range("H132").Value = 3475.5925
$ws.Range("I132").Value = 2163.5
$ws.Range("J132").Value = 4028.0527
$ws.Range("K132").Value = 6490.5
$ws.Range("L132").Value = 12084.1581
$ws.Range("M132").Value = -3960.5
$ws.Range("N132").Value = -17144.1581

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5170.5264
$ws.Range("I132").Value = 2065.5715
$ws.Range("J132").Value = 9006.058999999999
$ws.Range("K132").Value = 6196.7145
$ws.Range("L132").Value = 27018.177
$ws.Range("M132").Value = -3666.7145
$ws.Range("N132").Value = -32078.177
$ws.Range("H136").Value = 956.43396
$ws.Range("I136").Value = 652.0476
$ws.Range("K136").Value = 1956.1428
$ws.Range("M136").Value = 593.8571999999999
